$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 470, shifting existing rows 470-487 down to 471-488.
$ws.Rows.Item(470).Insert()

# Populate the newly inserted row 470 with the new weekly price record.
$ws.Cells.Item(470, 1).Value = 2
$ws.Cells.Item(470, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(470, 3).Value = "Coquimbo"
$ws.Cells.Item(470, 4).Value2 = 45267
$ws.Cells.Item(470, 5).Value = 4
$ws.Cells.Item(470, 6).Value = 100112021
$ws.Cells.Item(470, 7).Value = "Ají"
$ws.Cells.Item(470, 8).Value = "Americana (o)"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 160
$ws.Cells.Item(470, 11).Value = 28000
$ws.Cells.Item(470, 12).Value = 30000
$ws.Cells.Item(470, 13).Value = 29000
$ws.Cells.Item(470, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(470, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(470, 16).Value = 1160
$ws.Cells.Item(470, 17).Value = 25
$ws.Cells.Item(470, 18).Value = "Hortaliza"
